# Project proposal content rewrite: fill in the Goal / Background / Overview /
# Current State / Feedback Requested sections with real content instead of the
# placeholder prompt questions.
#
# Applied bottom-to-top so paragraph indices stay valid as we go (the Overview
# edit removes a paragraph, which would shift everything below it).

$d = $word.ActiveDocument
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Feedback Requested answer: plain text swap, no paragraph merge ---
[void]$d.Content.Find.Execute(
    'Identify the specific feedback requested. For example, is the goal to determine if the specified "why" something that resonates with users? Do we know if others have developed similar work? Are community contributions critical to the project?',
    $false, $false, $false, $false, $false, $true, 1, $false,
    'Do you have tools you would consider contributing?', 2)

# --- Current State answer: replace text and absorb the trailing blank " " paragraph ---
$pCS = $d.Paragraphs(11)
$pBlank = $d.Paragraphs(12)
$rngCS = $d.Range($pCS.Range.Start, $pBlank.Range.End)
$xmlCS = $pkgOpen + '<w:p><w:r><w:t>This project has very little progress. A few packages of tools have been developed for personal use, but the intent of this project is to share these tools.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' + $pkgClose
[void]$rngCS.InsertXML($xmlCS)

# --- Overview answer: replace text AND delete the trailing blank " " paragraph entirely ---
$pOV = $d.Paragraphs(8)
$pBlank2 = $d.Paragraphs(9)
$rngOV = $d.Range($pOV.Range.Start, $pBlank2.Range.End)
$xmlOV = $pkgOpen + '<w:p><w:r><w:t xml:space="preserve">At present, the project is in an exploratory stage. The concept is that individual contributors to this project should create a fork and add their own tools to the project in a way that is as organized as possible. Based on those forks and contributions, project maintainers will attempt to organize the </w:t></w:r><w:r><w:t xml:space="preserve">tools in a logical fashion and </w:t></w:r><w:r><w:t>then pull changes from other contributors into the main line of the project. These will then be packaged and made available for distribution.</w:t></w:r></w:p>' + $pkgClose
[void]$rngOV.InsertXML($xmlOV)

# --- Background answer: replace text (with spell-check markers around OpenG/GPower) and absorb the trailing blank " " paragraph ---
$pBG = $d.Paragraphs(5)
$pBlank3 = $d.Paragraphs(6)
$rngBG = $d.Range($pBG.Range.Start, $pBlank3.Range.End)
$xmlBG = $pkgOpen + '<w:p><w:r><w:t xml:space="preserve">Our organization does a reasonable job of promoting reuse of reference designs and specific APIs, but does not do a good job of sharing or aggregating less complete sets of tools. The LabVIEW community has done a good job of this in the past, with toolsets like </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>OpenG</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, MGI, and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GPower</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. This project does not aim to replace these tools, but to instead make additional toolsets available to the community.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' + $pkgClose
[void]$rngBG.InsertXML($xmlBG)

# --- Goal answer: replace text and absorb the trailing blank " " paragraph ---
$pG = $d.Paragraphs(2)
$pBlank4 = $d.Paragraphs(3)
$rngG = $d.Range($pG.Range.Start, $pBlank4.Range.End)
$xmlG = $pkgOpen + '<w:p><w:r><w:t>Develop a large set of palettes which provide common tools that can be shared across our organization and with our partners and customers.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' + $pkgClose
[void]$rngG.InsertXML($xmlG)
